$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.955.14'
$ws.Range("E2").Value = '  +7.05%  '

# Row 3
$ws.Range("D3").Value = '3.660.83'
$ws.Range("E3").Value = '  +19.23%  '

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").Value = "'599.60"
$ws.Range("E5").Value = '  +4.25%  '

# Row 6
$ws.Range("D6").Value = "'186.29"
$ws.Range("E6").Value = '  +9.24%  '

# Row 7
$ws.Range("D7").Value = '3.661.62'
$ws.Range("E7").Value = '  +19.39%  '

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = '  -0.12%  '

# Row 9
$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = '  +5.33%  '

# Row 10
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = '  +10.64%  '

# Row 11
$ws.Range("D11").Value = "'6.55"
$ws.Range("E11").Value = '  +4.32%  '

# Row 12
$ws.Range("D12").Value = "'0.499"
$ws.Range("E12").Value = '  +6.55%  '

# Row 13
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = "'39.69"
$ws.Range("E13").Value = '  +11.25%  '

# Row 14
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").Value = "'0.0000258"
$ws.Range("E14").Value = '  +8.04%  '

# Row 15
$ws.Range("D15").Value = '4.269.23'
$ws.Range("E15").Value = '  +19.16%  '

# Row 16
$ws.Range("D16").Value = '3.648.15'
$ws.Range("E16").Value = '  +18.70%  '

# Row 17
$ws.Range("D17").Value = '70.555.45'
$ws.Range("E17").Value = '  +6.53%  '

# Row 18
$ws.Range("E18").Value = '  +2.29%  '

# Row 19
$ws.Range("D19").Value = "'7.55"
$ws.Range("E19").Value = '  +8.85%  '

# Row 20
$ws.Range("D20").Value = "'17.42"
$ws.Range("E20").Value = '  +5.02%  '

# Row 21
$ws.Range("D21").Value = "'513.49"
$ws.Range("E21").Value = '  +6.00%  '

# Row 22
$ws.Range("D22").Value = "'9.26"
$ws.Range("E22").Value = '  +21.24%  '

# Row 23
$ws.Range("D23").Value = "'0.753"
$ws.Range("E23").Value = '  +10.09%  '

# Row 24
$ws.Range("D24").Value = "'88.65"
$ws.Range("E24").Value = '  +7.81%  '

# Row 25
$ws.Range("D25").Value = "'13.61"
$ws.Range("E25").Value = '  +8.09%  '

# Row 26
$ws.Range("D26").Value = "'2.42"
$ws.Range("E26").Value = '  +10.28%  '

# Row 27
$ws.Range("D27").Value = "'10.85"
$ws.Range("E27").Value = '  +8.09%  '

# Row 28
$ws.Range("E28").Value = '  +0.14%  '

# Row 29
$ws.Range("D29").Value = "'2.57"
$ws.Range("E29").Value = '  +14.72%  '

# Row 30
$ws.Range("D30").Value = "'8.30"
$ws.Range("E30").Value = '  +6.22%  '

# Row 31
$ws.Range("D31").Value = "'32.19"
$ws.Range("E31").Value = '  +16.55%  '

# Row 32
$ws.Range("D32").Value = "'2.78"
$ws.Range("E32").Value = '  +7.50%  '

# Row 33
$ws.Range("D33").Value = "'0.0000110"
$ws.Range("E33").Value = '  +20.13%  '

# Row 34
$ws.Range("E34").Value = '  +6.34%  '

# Row 35
$ws.Range("D35").Value = "'0.994"
$ws.Range("E35").Value = '  -0.51%  '

# Row 36
$ws.Range("D36").Value = "'6.19"
$ws.Range("E36").Value = '  +11.65%  '

# Row 37
$ws.Range("D37").Value = "'1.03"
$ws.Range("E37").Value = '  +9.24%  '

# Row 38
$ws.Range("D38").Value = "'0.337"
$ws.Range("E38").Value = '  +12.36%  '

# Row 39
$ws.Range("D39").Value = "'2.13"
$ws.Range("E39").Value = '  +8.90%  '

# Row 40
$ws.Range("D40").Value = "'47.41"
$ws.Range("E40").Value = '  -0.66%  '

# Row 41
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = "'0.129"
$ws.Range("E41").Value = '  +5.91%  '

# Row 42
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = "'50.98"
$ws.Range("E42").Value = '  +3.94%  '

# Row 43
$ws.Range("D43").Value = "'8.95"
$ws.Range("E43").Value = '  +8.95%  '

# Row 44
$ws.Range("D44").Value = '3.153.41'
$ws.Range("E44").Value = '  +13.65%  '

# Row 45
$ws.Range("D45").Value = "'2.81"
$ws.Range("E45").Value = '  +11.36%  '

# Row 46
$ws.Range("D46").Value = "'408.65"
$ws.Range("E46").Value = '  +12.37%  '

# Row 47
$ws.Range("D47").Value = "'0.0367"
$ws.Range("E47").Value = '  +7.05%  '

# Row 48
$ws.Range("D48").Value = "'28.03"
$ws.Range("E48").Value = '  +15.96%  '

# Row 49
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").Value = "'2.49"
$ws.Range("E49").Value = '  +16.26%  '

# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = "'134.69"
$ws.Range("E50").Value = '  +0.06%  '

# Row 51
$ws.Range("E51").Value = '  +0.03%  '
